$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-22 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-23 Wednesday", 2)
$d.Content.Find.Execute("45×93=", $true, $false, $false, $false, $false, $true, 1, $false, "87×52=", 2)
$d.Content.Find.Execute("39×79=", $true, $false, $false, $false, $false, $true, 1, $false, "85×87=", 2)
$d.Content.Find.Execute("34×80=", $true, $false, $false, $false, $false, $true, 1, $false, "20×94=", 2)
$d.Content.Find.Execute("27×85=", $true, $false, $false, $false, $false, $true, 1, $false, "98×79=", 2)
$d.Content.Find.Execute("96×66=", $true, $false, $false, $false, $false, $true, 1, $false, "30×28=", 2)
$d.Content.Find.Execute("45×23=", $true, $false, $false, $false, $false, $true, 1, $false, "83×14=", 2)
$d.Content.Find.Execute("79×89=", $true, $false, $false, $false, $false, $true, 1, $false, "73×56=", 2)
$d.Content.Find.Execute("60×98=", $true, $false, $false, $false, $false, $true, 1, $false, "63×14=", 2)
$d.Content.Find.Execute("79×90=", $true, $false, $false, $false, $false, $true, 1, $false, "18×94=", 2)
$d.Content.Find.Execute("61×45=", $true, $false, $false, $false, $false, $true, 1, $false, "89×45=", 2)
$d.Content.Find.Execute("55×35=", $true, $false, $false, $false, $false, $true, 1, $false, "13×17=", 2)
$d.Content.Find.Execute("51×68=", $true, $false, $false, $false, $false, $true, 1, $false, "36×70=", 2)
$d.Content.Find.Execute("88×24=", $true, $false, $false, $false, $false, $true, 1, $false, "12×37=", 2)
$d.Content.Find.Execute("79×16=", $true, $false, $false, $false, $false, $true, 1, $false, "87×52=", 2)
$d.Content.Find.Execute("33×18=", $true, $false, $false, $false, $false, $true, 1, $false, "31×22=", 2)
$d.Content.Find.Execute("66×63=", $true, $false, $false, $false, $false, $true, 1, $false, "36×54=", 2)
$d.Content.Find.Execute("12×70=", $true, $false, $false, $false, $false, $true, 1, $false, "82×80=", 2)
$d.Content.Find.Execute("29×87=", $true, $false, $false, $false, $false, $true, 1, $false, "73×81=", 2)
$d.Content.Find.Execute("77×56=", $true, $false, $false, $false, $false, $true, 1, $false, "80×19=", 2)
$d.Content.Find.Execute("25×15=", $true, $false, $false, $false, $false, $true, 1, $false, "55×98=", 2)
$d.Content.Find.Execute("17×60=", $true, $false, $false, $false, $false, $true, 1, $false, "50×59=", 2)
$d.Content.Find.Execute("41×84=", $true, $false, $false, $false, $false, $true, 1, $false, "36×84=", 2)
$d.Content.Find.Execute("58×71=", $true, $false, $false, $false, $false, $true, 1, $false, "50×11=", 2)
$d.Content.Find.Execute("72×32=", $true, $false, $false, $false, $false, $true, 1, $false, "93×16=", 2)
$d.Content.Find.Execute("78×84=", $true, $false, $false, $false, $false, $true, 1, $false, "77×37=", 2)
